$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns at E, G, I (original E/F/G/H/I shift right to F/H/J/K/L)
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(9).Insert()

# Header row (row 1) - fill the three newly inserted header cells
$ws.Range("E1").Value = "Answer w/ Explanation"
$ws.Range("G1").Value = "Time w/ Explanation"
$ws.Range("I1").Value = "Confidence w/ Explanation"

# Row 2 - new "n/a" cells in inserted columns
$ws.Range("E2").Value = "n/a"
$ws.Range("G2").Value = "n/a"
$ws.Range("I2").Value = "n/a"
$ws.Range("K2").Value = 1

# Row 3 - new "n/a" cells in inserted columns
$ws.Range("E3").Value = "n/a"
$ws.Range("G3").Value = "n/a"
$ws.Range("I3").Value = "n/a"
$ws.Range("K3").Value = 1

# Row 4 - new "n/a" cells in inserted columns
$ws.Range("E4").Value = "n/a"
$ws.Range("G4").Value = "n/a"
$ws.Range("I4").Value = "n/a"
$ws.Range("K4").Value = 0

# Row 6 - new participant data
$ws.Range("B6").Value = "mission10"
$ws.Range("C6").Value = "agentPolicy0"
$ws.Range("D6").Value = "no"
$ws.Range("E6").Value = "yes"
$ws.Range("F6").Value = 3.15
$ws.Range("G6").Value = 3.06
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "yes"
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = "Think agentPolicy0 is good, but maybe there is another better policy. Wonder what happens if robot goes faster from L1-L2. Explanation changed mind."

# Row 7 - new participant data
$ws.Range("B7").Value = "mission20"
$ws.Range("C7").Value = "agentPolicy1"
$ws.Range("D7").Value = "yes"
$ws.Range("E7").Value = "yes"
$ws.Range("F7").Value = 1
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = "yes"
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = "Quick to answer: intrusive penalty is very high especially compared to time cost. Explanation didn't change mind; just confirmed."

# Row 8 - new participant data
$ws.Range("B8").Value = "mission30"
$ws.Range("C8").Value = "agentPolicy2"
$ws.Range("D8").Value = "no"
$ws.Range("E8").Value = "no"
$ws.Range("F8").Value = 2.23
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = "no"
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = "Explanation didn't change mind; just confirmed."

# Row 9 - new participant data
$ws.Range("B9").Value = "mission40"
$ws.Range("C9").Value = "agentPolicy3"
$ws.Range("D9").Value = "no"
$ws.Range("E9").Value = "no"
$ws.Range("F9").Value = 2.04
$ws.Range("G9").Value = 3.48
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = "no"
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = "Think another policy is a better but unsure; think they are close. Explanation improved confidence."

# Match final selection state from the diff
[void]$ws.Range("L6").Select()
